# Brisbane_stats.xlsx update:
#   The sheet holds one row per statistic and one column per match, with the
#   most-recently-played match duplicated into the final column each time the
#   dataset is refreshed (so JO == old JP going into this edit). This commit
#   fixes a finals-round numbering mistake and streamlines the pipeline so the
#   "current match" value only ever appears once more at the very end -
#   structurally that shows up as three extra columns (JP..JR) being inserted
#   right before the old trailing column (which shifts from JP to JS), all
#   four of them carrying the same value as JO.
#
# Net effect required:
#   - dimension grows from A1:JP102 to A1:JS102
#   - for every row r (1..102): JP, JQ, JR all get the same value that JO/JS
#     already hold (JS is just the old JP, shifted three columns right)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataCol = 275   # JO - source of the value being propagated
$insertAt     = 276   # JP - where the 3 new columns get inserted
$insertCount  = 3
$lastRow      = 102

# Insert 3 new columns at JP, shifting JP:JP102 (and everything after it)
# three columns to the right - JP becomes JS, dimension becomes A1:JS102.
$ws.Range("JP1:JR1").EntireColumn.Insert(-4161)

# The 3 freshly inserted columns (now JP, JQ, JR) are blank; fill each one,
# row by row, with the value already sitting in JO (which is also what the
# shifted-over JS column carries).
for ($r = 1; $r -le $lastRow; $r++) {
  $v = $ws.Cells.Item($r, $firstDataCol).Value2
  for ($i = 0; $i -lt $insertCount; $i++) {
    $ws.Cells.Item($r, $insertAt + $i).Value2 = $v
  }
}
